# Update marksheet "Right" and "Max" figures on the quiz worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 "Marking" - Right column (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Row 12 "Total" - Right column (B12): 72 -> 120
$ws.Range("B12").Value = 120

# Row 12 "Total" - Max column (E12): "70/84" -> "120/140"
$ws.Range("E12").Value = "120/140"
